# Updated cryptos list on Sat Nov  4 01:37:04 UTC 2023 with GitHub Actions
#
# Applies the per-cell Coin/Link/Price/Volume(1h) updates described by the
# upstream diff. Price cells in column D are sometimes numeric-looking
# strings (e.g. "231.35", "0.0690") that must remain literal text (matching
# the workbook's inlineStr cells) rather than being auto-coerced to numbers
# by Excel's normal cell-value parsing, so those are written via a
# Text-number-format round trip and then have their formatting cleared back
# to the sheet default.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '34.940.68'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +0.40%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.839.68'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  +2.16%  '
$ws.Range('E4').Value = '  +0.34%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '231.35'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.31%  '
$ws.Range('E6').Value = '  +2.30%  '
$ws.Range('E7').Value = '  +0.23%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '40.29'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +1.84%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.329'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +2.35%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0690'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +2.10%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0981'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -0.76%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '2.108.20'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +2.34%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '11.44'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +5.66%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '1.844.48'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +2.44%  '
$ws.Range('E15').Value = '  +2.61%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '4.64'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +1.93%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '34.967.00'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +0.84%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '69.76'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +1.43%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.0₃0790'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +1.90%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '239.49'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +1.42%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '12.17'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +4.04%  '
$ws.Range('E22').Value = '  +1.71%  '
$ws.Range('E23').Value = '  +0.14%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.28'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +2.97%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '171.99'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +0.84%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '7.79'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +1.33%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '17.48'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +1.94%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.123'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +3.85%  '
$ws.Range('E29').Value = '  +0.26%  '
$ws.Range('E30').Value = '  +0.42%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.0552'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +1.23%  '
$ws.Range('E32').Value = '  -1.30%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.94'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +1.14%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.95'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +12.11%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.58'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +22.99%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.758'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +12.84%  '
$ws.Range('B37').Value = 'TrustWalletToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.22'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +8.07%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '1.07'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +11.96%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '89.28'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -0.71%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.348.85'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +3.59%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.0195'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +2.53%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '14.62'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +2.21%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '2.27'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +3.86%  '
$ws.Range('E44').Value = '  -1.96%  '
$ws.Range('E45').Value = '  +0.66%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0529'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +3.97%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '6.24'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +1.79%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.025.77'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +2.50%  '
$ws.Range('B49').Value = 'PaxDollar'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.01'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +0.19%  '
$ws.Range('B50').Value = 'THORChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '3.40'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +18.81%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.0668'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +1.58%  '
